# Objective function cost definition change:
# The "variability" sheet's A2 driver value changes from 0.25 to 0.75.
# All dependent formulas (1-variability!$A$2 / 1+variability!$A$2) used on
# scenarios_base, three_scenarios and three_scenarios_new recalculate
# automatically from this single change.

$wb = $excel.ActiveWorkbook

$variabilitySheet    = $wb.Worksheets.Item("variability")
$newScenariosSheet   = $wb.Worksheets.Item("three_scenarios_new")

$variabilitySheet.Range("A2").Value = 0.75

# Update the saved selection on "three_scenarios_new" (was J8, now E2) while
# it is no longer the active tab.
$newScenariosSheet.Activate()
$newScenariosSheet.Range("E2").Select()

# "variability" becomes the active/selected sheet (selection stays at A3).
$variabilitySheet.Activate()
$variabilitySheet.Range("A3").Select()
